$d = $word.ActiveDocument

# Append two new paragraphs after the last paragraph of the document,
# matching the formatting (pt-BR language run properties) of the
# existing paragraphs.

$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$p1 = $d.Paragraphs.Last
$p1.Range.Text = "*Criar mensagens de erro nos inputs do form"

$p1.Range.InsertParagraphAfter()

$p2 = $d.Paragraphs.Last
$p2.Range.Text = "*Desabilitar botão até que o form seja válido"
